# debug koperasi and import
# The "nama" header label (A1) was relabeled to "name". All other rows
# (T1..T5 Items) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
